$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Updated at" footer timestamp: 12:52 -> 13:22
$ws.Range("A1").Value = 'Datos actualizados a 13 de Abril de 2020 a las 13:22'

# Refreshed COVID-19 per-province counters (Casos totales / activos / Recuperados / Muertes).
# Several provinces changed rank order as totals moved, so some rows also get a new
# "Ciudad" label (the underlying numbers stay attached to their own data row).

# Row 7: Bizkaia/Vizcaya
$ws.Range("B7").Value = 5933
$ws.Range("C7").Value = 5026
$ws.Range("D7").Value = 5161
$ws.Range("E7").Value = 431

# Row 9: Valencia/Valencia
$ws.Range("B9").Value = 4727
$ws.Range("C9").Value = 1487
$ws.Range("D9").Value = 2815
$ws.Range("E9").Value = 425

# Row 13: was Zaragoza, now Alacant/Alicante (ranking swap with row 14)
$ws.Range("A13").Value = 'Alacant/Alicante'
$ws.Range("B13").Value = 3212
$ws.Range("C13").Value = 1046
$ws.Range("D13").Value = 1824
$ws.Range("E13").Value = 342

# Row 14: was Alacant/Alicante, now Zaragoza
$ws.Range("A14").Value = 'Zaragoza'
$ws.Range("B14").Value = 3137
$ws.Range("C14").Value = 717
$ws.Range("D14").Value = 2075
$ws.Range("E14").Value = 345

# Row 16: Araba/Alava
$ws.Range("B16").Value = 3034
$ws.Range("C16").Value = 5026
$ws.Range("D16").Value = 5161
$ws.Range("E16").Value = 266

# Row 21: was Sevilla, now Gipuzkoa/Guipuzcoa
$ws.Range("A21").Value = 'Gipuzkoa/Guipuzcoa'
$ws.Range("B21").Value = 2051
$ws.Range("C21").Value = 5026
$ws.Range("D21").Value = 5161
$ws.Range("E21").Value = 134

# Row 22: was A Coruña, now Sevilla
$ws.Range("A22").Value = 'Sevilla'
$ws.Range("B22").Value = 2034
$ws.Range("C22").Value = 247
$ws.Range("D22").Value = 1618
$ws.Range("E22").Value = 169

# Row 23: was Asturias, now A Coruña
$ws.Range("A23").Value = 'A Coruña'
$ws.Range("B23").Value = 1969
$ws.Range("C23").Value = 333
$ws.Range("D23").Value = 1788
$ws.Range("E23").Value = 67

# Row 24: was Gipuzkoa/Guipuzcoa, now Asturias
$ws.Range("A24").Value = 'Asturias'
$ws.Range("B24").Value = 1958
$ws.Range("C24").Value = 434
$ws.Range("D24").Value = 1375
$ws.Range("E24").Value = 149

# Row 27: was Granada, now Caceres
$ws.Range("A27").Value = 'Caceres'
$ws.Range("B27").Value = 1776
$ws.Range("C27").Value = 237
$ws.Range("D27").Value = 1276
$ws.Range("E27").Value = 263

# Row 28: was Caceres, now Granada
$ws.Range("A28").Value = 'Granada'
$ws.Range("B28").Value = 1772
$ws.Range("C28").Value = 317
$ws.Range("D28").Value = 1290
$ws.Range("E28").Value = 165

# Row 37: Guadalajara
$ws.Range("B37").Value = 1116
$ws.Range("C37").Value = 270
$ws.Range("D37").Value = 737
$ws.Range("E37").Value = 109

# Row 43: Badajoz
$ws.Range("B43").Value = 882
$ws.Range("C43").Value = 274
$ws.Range("D43").Value = 550
$ws.Range("E43").Value = 58
